$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove existing merged cells so row/cell edits are unobstructed
$ws.Cells.UnMerge()

# 2. Update cell values for rows 12-55 to match the reorganised specification rows
# Row 12
$ws.Range("A12").ClearContents()
$ws.Range("B12").ClearContents()
$ws.Range("C12").Value = 'Application'
$ws.Range("D12").Value = 'Documents[]'
$ws.Range("E12").Value = 'File'
$ws.Range("F12").Value = 'Base64'
$ws.Range("G12").Value = 'Base64-encoded content of the file for inline file uploads'
$ws.Range("H12").Value = 'string'
$ws.Range("I12").Value = 'MAY'

# Row 13
$ws.Range("A13").ClearContents()
$ws.Range("B13").ClearContents()
$ws.Range("C13").Value = 'Application'
$ws.Range("D13").Value = 'Documents[]'
$ws.Range("E13").Value = 'File'
$ws.Range("F13").Value = 'Filename'
$ws.Range("G13").Value = 'Name of the file being uploaded'
$ws.Range("H13").Value = 'string'
$ws.Range("I13").Value = 'MUST'

# Row 14
$ws.Range("A14").ClearContents()
$ws.Range("B14").ClearContents()
$ws.Range("C14").Value = 'Application'
$ws.Range("D14").Value = 'Documents[]'
$ws.Range("E14").Value = 'File'
$ws.Range("F14").Value = 'MIME type'
$ws.Range("G14").Value = 'The file''s MIME type such as application/pdf or image/jpeg'
$ws.Range("H14").Value = 'string'
$ws.Range("I14").Value = 'MAY'

# Row 15
$ws.Range("A15").ClearContents()
$ws.Range("B15").ClearContents()
$ws.Range("C15").Value = 'Application'
$ws.Range("D15").Value = 'Documents[]'
$ws.Range("E15").Value = 'File'
$ws.Range("F15").Value = 'Checksum'
$ws.Range("G15").Value = 'Hash of the file contents used for file validation and checking files have not been tampered with'
$ws.Range("H15").Value = 'string'
$ws.Range("I15").Value = 'MAY'

# Row 16
$ws.Range("A16").ClearContents()
$ws.Range("B16").ClearContents()
$ws.Range("C16").Value = 'Application'
$ws.Range("D16").Value = 'Documents[]'
$ws.Range("E16").Value = 'File'
$ws.Range("F16").Value = 'File size'
$ws.Range("G16").Value = 'Size of the file in bytes that can be used to enforce limits'
$ws.Range("H16").Value = 'number'
$ws.Range("I16").Value = 'MAY'

# Row 17
$ws.Range("A17").ClearContents()
$ws.Range("B17").ClearContents()
$ws.Range("C17").Value = 'Application'
$ws.Range("D17").Value = 'Fee'
$ws.Range("E17").Value = 'Amount'
$ws.Range("F17").ClearContents()
$ws.Range("G17").Value = 'The total amount due for the application fee'
$ws.Range("H17").Value = 'number'
$ws.Range("I17").Value = 'MUST'

# Row 18
$ws.Range("A18").ClearContents()
$ws.Range("B18").ClearContents()
$ws.Range("C18").Value = 'Application'
$ws.Range("D18").Value = 'Fee'
$ws.Range("E18").Value = 'Amount paid'
$ws.Range("F18").ClearContents()
$ws.Range("G18").Value = 'The amount paid towards the application fee'
$ws.Range("H18").Value = 'number'
$ws.Range("I18").Value = 'MUST'

# Row 19
$ws.Range("A19").ClearContents()
$ws.Range("B19").ClearContents()
$ws.Range("C19").Value = 'Application'
$ws.Range("D19").Value = 'Fee'
$ws.Range("E19").Value = 'Transactions[]'
$ws.Range("F19").ClearContents()
$ws.Range("G19").Value = 'References to payments or financial transactions related to this application'
$ws.Range("H19").Value = 'string'
$ws.Range("I19").Value = 'MAY'

# Row 20
$ws.Range("A20").Value = 'Agent contact details'
$ws.Range("B20").Value = 'Name and contact information if an agent is being used.'
$ws.Range("C20").Value = 'Agent reference'
$ws.Range("D20").ClearContents()
$ws.Range("E20").ClearContents()
$ws.Range("F20").ClearContents()
$ws.Range("G20").Value = 'A reference to an agent object'
$ws.Range("H20").Value = 'string'
$ws.Range("I20").Value = 'MUST'

# Row 21
$ws.Range("A21").ClearContents()
$ws.Range("B21").ClearContents()
$ws.Range("C21").Value = 'Contact details'
$ws.Range("D21").Value = 'Email'
$ws.Range("E21").ClearContents()
$ws.Range("F21").ClearContents()
$ws.Range("G21").Value = 'The email address that can be used for electronic correspondence with the individual'
$ws.Range("H21").Value = 'string'
$ws.Range("I21").Value = 'MUST'

# Row 22
$ws.Range("A22").ClearContents()
$ws.Range("B22").ClearContents()
$ws.Range("C22").Value = 'Contact details'
$ws.Range("D22").Value = 'Phone number(s)[]'
$ws.Range("E22").Value = 'Phone number'
$ws.Range("F22").ClearContents()
$ws.Range("G22").Value = 'A phone number'
$ws.Range("H22").Value = 'string'
$ws.Range("I22").Value = 'MAY'

# Row 23
$ws.Range("A23").ClearContents()
$ws.Range("B23").ClearContents()
$ws.Range("C23").Value = 'Contact details'
$ws.Range("D23").Value = 'Phone number(s)[]'
$ws.Range("E23").Value = 'Contact priority'
$ws.Range("F23").ClearContents()
$ws.Range("G23").Value = 'The priority of a number'
$ws.Range("H23").Value = 'enum'
$ws.Range("I23").Value = 'MAY'

# Row 24
$ws.Range("A24").Value = 'Agent details'
$ws.Range("B24").Value = 'Name and contact information if an agent is being used.'
$ws.Range("C24").Value = 'agent'
$ws.Range("D24").Value = 'Reference'
$ws.Range("E24").ClearContents()
$ws.Range("F24").ClearContents()
$ws.Range("G24").Value = 'A unique reference for the data item'
$ws.Range("H24").Value = 'string'
$ws.Range("I24").Value = 'MUST'

# Row 25
$ws.Range("A25").ClearContents()
$ws.Range("B25").ClearContents()
$ws.Range("C25").Value = 'agent'
$ws.Range("D25").Value = 'Person'
$ws.Range("E25").Value = 'Title'
$ws.Range("F25").ClearContents()
$ws.Range("G25").Value = 'The title of the individual'
$ws.Range("H25").Value = 'string'
$ws.Range("I25").Value = 'MAY'

# Row 26
$ws.Range("A26").ClearContents()
$ws.Range("B26").ClearContents()
$ws.Range("C26").Value = 'agent'
$ws.Range("D26").Value = 'Person'
$ws.Range("E26").Value = 'First Name'
$ws.Range("F26").ClearContents()
$ws.Range("G26").Value = 'The first name of the individual'
$ws.Range("H26").Value = 'string'
$ws.Range("I26").Value = 'MUST'

# Row 27
$ws.Range("A27").ClearContents()
$ws.Range("B27").ClearContents()
$ws.Range("C27").Value = 'agent'
$ws.Range("D27").Value = 'Person'
$ws.Range("E27").Value = 'Last Name'
$ws.Range("F27").ClearContents()
$ws.Range("G27").Value = 'The last name of the individual'
$ws.Range("H27").Value = 'string'
$ws.Range("I27").Value = 'MUST'

# Row 28
$ws.Range("A28").ClearContents()
$ws.Range("B28").ClearContents()
$ws.Range("C28").Value = 'agent'
$ws.Range("D28").Value = 'Person'
$ws.Range("E28").Value = 'Address Text'
$ws.Range("F28").ClearContents()
$ws.Range("G28").Value = 'Flexible field for capturing addresses'
$ws.Range("H28").Value = 'string'
$ws.Range("I28").Value = 'MUST'

# Row 29
$ws.Range("A29").ClearContents()
$ws.Range("B29").ClearContents()
$ws.Range("C29").Value = 'agent'
$ws.Range("D29").Value = 'Person'
$ws.Range("E29").Value = 'Postcode'
$ws.Range("F29").ClearContents()
$ws.Range("G29").Value = 'The postal code'
$ws.Range("H29").Value = 'string'
$ws.Range("I29").Value = 'MAY'

# Row 30
$ws.Range("A30").ClearContents()
$ws.Range("B30").ClearContents()
$ws.Range("C30").Value = 'agent'
$ws.Range("D30").Value = 'Company'
$ws.Range("E30").ClearContents()
$ws.Range("F30").ClearContents()
$ws.Range("G30").Value = 'The name of a company (that the agent works for)'
$ws.Range("H30").Value = 'string'
$ws.Range("I30").Value = 'MAY'

# Row 31
$ws.Range("A31").ClearContents()
$ws.Range("B31").ClearContents()
$ws.Range("C31").Value = 'agent'
$ws.Range("D31").Value = 'User role'
$ws.Range("E31").ClearContents()
$ws.Range("F31").ClearContents()
$ws.Range("G31").Value = 'The role of the named individual. Agent or proxy'
$ws.Range("H31").Value = 'enum'
$ws.Range("I31").Value = 'MAY'

# Row 32
$ws.Range("A32").Value = 'Applicant contact details'
$ws.Range("B32").Value = 'Telephone number and email address of the applicant.'
$ws.Range("C32").Value = 'Applicant reference'
$ws.Range("D32").ClearContents()
$ws.Range("E32").ClearContents()
$ws.Range("F32").ClearContents()
$ws.Range("G32").Value = 'Reference to match contact details to a named individual from the applicant details component'
$ws.Range("H32").Value = 'string'
$ws.Range("I32").Value = 'MUST'

# Row 33
$ws.Range("A33").ClearContents()
$ws.Range("B33").ClearContents()
$ws.Range("C33").Value = 'Contact details'
$ws.Range("D33").Value = 'Email'
$ws.Range("E33").ClearContents()
$ws.Range("F33").ClearContents()
$ws.Range("G33").Value = 'The email address that can be used for electronic correspondence with the individual'
$ws.Range("H33").Value = 'string'
$ws.Range("I33").Value = 'MUST'

# Row 34
$ws.Range("A34").ClearContents()
$ws.Range("B34").ClearContents()
$ws.Range("C34").Value = 'Contact details'
$ws.Range("D34").Value = 'Phone number(s)[]'
$ws.Range("E34").Value = 'Phone number'
$ws.Range("F34").ClearContents()
$ws.Range("G34").Value = 'A phone number'
$ws.Range("H34").Value = 'string'
$ws.Range("I34").Value = 'MAY'

# Row 35
$ws.Range("A35").ClearContents()
$ws.Range("B35").ClearContents()
$ws.Range("C35").Value = 'Contact details'
$ws.Range("D35").Value = 'Phone number(s)[]'
$ws.Range("E35").Value = 'Contact priority'
$ws.Range("F35").ClearContents()
$ws.Range("G35").Value = 'The priority of a number'
$ws.Range("H35").Value = 'enum'
$ws.Range("I35").Value = 'MAY'

# Row 36
$ws.Range("A36").Value = 'Applicant details'
$ws.Range("B36").Value = 'Name and contact information for the parties making the application.'
$ws.Range("C36").Value = 'Applicants[]'
$ws.Range("D36").Value = 'Reference'
$ws.Range("E36").ClearContents()
$ws.Range("F36").ClearContents()
$ws.Range("G36").Value = 'A unique reference for the data item'
$ws.Range("H36").Value = 'string'
$ws.Range("I36").Value = 'MUST'

# Row 37
$ws.Range("A37").ClearContents()
$ws.Range("B37").ClearContents()
$ws.Range("C37").Value = 'Applicants[]'
$ws.Range("D37").Value = 'Person'
$ws.Range("E37").Value = 'Title'
$ws.Range("F37").ClearContents()
$ws.Range("G37").Value = 'The title of the individual'
$ws.Range("H37").Value = 'string'
$ws.Range("I37").Value = 'MAY'

# Row 38
$ws.Range("A38").ClearContents()
$ws.Range("B38").ClearContents()
$ws.Range("C38").Value = 'Applicants[]'
$ws.Range("D38").Value = 'Person'
$ws.Range("E38").Value = 'First Name'
$ws.Range("F38").ClearContents()
$ws.Range("G38").Value = 'The first name of the individual'
$ws.Range("H38").Value = 'string'
$ws.Range("I38").Value = 'MUST'

# Row 39
$ws.Range("A39").ClearContents()
$ws.Range("B39").ClearContents()
$ws.Range("C39").Value = 'Applicants[]'
$ws.Range("D39").Value = 'Person'
$ws.Range("E39").Value = 'Last Name'
$ws.Range("F39").ClearContents()
$ws.Range("G39").Value = 'The last name of the individual'
$ws.Range("H39").Value = 'string'
$ws.Range("I39").Value = 'MUST'

# Row 40
$ws.Range("A40").ClearContents()
$ws.Range("B40").ClearContents()
$ws.Range("C40").Value = 'Applicants[]'
$ws.Range("D40").Value = 'Person'
$ws.Range("E40").Value = 'Address Text'
$ws.Range("F40").ClearContents()
$ws.Range("G40").Value = 'Flexible field for capturing addresses'
$ws.Range("H40").Value = 'string'
$ws.Range("I40").Value = 'MUST'

# Row 41
$ws.Range("A41").ClearContents()
$ws.Range("B41").ClearContents()
$ws.Range("C41").Value = 'Applicants[]'
$ws.Range("D41").Value = 'Person'
$ws.Range("E41").Value = 'Postcode'
$ws.Range("F41").ClearContents()
$ws.Range("G41").Value = 'The postal code'
$ws.Range("H41").Value = 'string'
$ws.Range("I41").Value = 'MAY'

# Row 42
$ws.Range("A42").Value = 'Conflict of interest'
$ws.Range("B42").Value = 'Details of any conflict of interest that may exist between the applicant and planning authority.'
$ws.Range("C42").Value = ''
$ws.Range("D42").Value = ''
$ws.Range("E42").Value = ''
$ws.Range("F42").Value = ''
$ws.Range("G42").Value = ''
$ws.Range("H42").Value = ''
$ws.Range("I42").Value = ''

# Row 43
$ws.Range("A43").Value = 'Checklist'
$ws.Range("B43").Value = 'Checking whether all the requirements of the form have been met, such as proof of payment or supporting documentation.'
$ws.Range("C43").Value = 'National requirement types[]'
$ws.Range("D43").ClearContents()
$ws.Range("E43").ClearContents()
$ws.Range("F43").ClearContents()
$ws.Range("G43").Value = 'List of the document types required for the given application type'
$ws.Range("H43").Value = 'string'
$ws.Range("I43").Value = 'MUST'

# Row 44
$ws.Range("A44").Value = 'Declaration'
$ws.Range("B44").Value = 'Signed and dated verification of the application''s accuracy.'
$ws.Range("C44").Value = 'Name'
$ws.Range("D44").ClearContents()
$ws.Range("E44").ClearContents()
$ws.Range("F44").ClearContents()
$ws.Range("G44").Value = 'A name of a person'
$ws.Range("H44").Value = 'string'
$ws.Range("I44").Value = 'MUST'

# Row 45
$ws.Range("A45").ClearContents()
$ws.Range("B45").ClearContents()
$ws.Range("C45").Value = 'Declaration confirmed'
$ws.Range("D45").ClearContents()
$ws.Range("E45").ClearContents()
$ws.Range("F45").ClearContents()
$ws.Range("G45").Value = 'Confirms the applicant or agent has reviewed and validated the information provided in the application'
$ws.Range("H45").Value = 'boolean'
$ws.Range("I45").Value = 'MUST'

# Row 46
$ws.Range("A46").ClearContents()
$ws.Range("B46").ClearContents()
$ws.Range("C46").Value = 'Declaration date'
$ws.Range("D46").ClearContents()
$ws.Range("E46").ClearContents()
$ws.Range("F46").ClearContents()
$ws.Range("G46").Value = 'The date the declaration was made'
$ws.Range("H46").Value = 'string'
$ws.Range("I46").Value = 'MUST'

# Row 47
$ws.Range("A47").Value = 'Site details'
$ws.Range("B47").Value = 'Where the proposed development will be built.'
$ws.Range("C47").Value = 'Site locations[]'
$ws.Range("D47").Value = 'Site boundary'
$ws.Range("E47").ClearContents()
$ws.Range("F47").ClearContents()
$ws.Range("G47").Value = 'Geometry of the site of the development, typically in GeoJSON format'
$ws.Range("H47").Value = 'wkt'
$ws.Range("I47").Value = 'MAY'

# Row 48
$ws.Range("A48").ClearContents()
$ws.Range("B48").ClearContents()
$ws.Range("C48").Value = 'Site locations[]'
$ws.Range("D48").Value = 'Address Text'
$ws.Range("E48").ClearContents()
$ws.Range("F48").ClearContents()
$ws.Range("G48").Value = 'Flexible field for capturing addresses'
$ws.Range("H48").Value = 'string'
$ws.Range("I48").Value = 'MAY'

# Row 49
$ws.Range("A49").ClearContents()
$ws.Range("B49").ClearContents()
$ws.Range("C49").Value = 'Site locations[]'
$ws.Range("D49").Value = 'Postcode'
$ws.Range("E49").ClearContents()
$ws.Range("F49").ClearContents()
$ws.Range("G49").Value = 'The postal code'
$ws.Range("H49").Value = 'string'
$ws.Range("I49").Value = 'MAY'

# Row 50
$ws.Range("A50").ClearContents()
$ws.Range("B50").ClearContents()
$ws.Range("C50").Value = 'Site locations[]'
$ws.Range("D50").Value = 'Easting'
$ws.Range("E50").ClearContents()
$ws.Range("F50").ClearContents()
$ws.Range("G50").Value = 'Easting coordinate in British National Grid (EPSG:27700)'
$ws.Range("H50").Value = 'number'
$ws.Range("I50").Value = 'MAY'

# Row 51
$ws.Range("A51").ClearContents()
$ws.Range("B51").ClearContents()
$ws.Range("C51").Value = 'Site locations[]'
$ws.Range("D51").Value = 'Northing'
$ws.Range("E51").ClearContents()
$ws.Range("F51").ClearContents()
$ws.Range("G51").Value = 'Northing coordinate in British National Grid (EPSG:27700)'
$ws.Range("H51").Value = 'number'
$ws.Range("I51").Value = 'MAY'

# Row 52
$ws.Range("A52").ClearContents()
$ws.Range("B52").ClearContents()
$ws.Range("C52").Value = 'Site locations[]'
$ws.Range("D52").Value = 'Latitude'
$ws.Range("E52").ClearContents()
$ws.Range("F52").ClearContents()
$ws.Range("G52").Value = 'Latitude coordinate in WGS84 (EPSG:4326)'
$ws.Range("H52").Value = 'number'
$ws.Range("I52").Value = 'MAY'

# Row 53
$ws.Range("A53").ClearContents()
$ws.Range("B53").ClearContents()
$ws.Range("C53").Value = 'Site locations[]'
$ws.Range("D53").Value = 'Longitude'
$ws.Range("E53").ClearContents()
$ws.Range("F53").ClearContents()
$ws.Range("G53").Value = 'Longitude coordinate in WGS84 (EPSG:4326)'
$ws.Range("H53").Value = 'number'
$ws.Range("I53").Value = 'MAY'

# Row 54
$ws.Range("A54").ClearContents()
$ws.Range("B54").ClearContents()
$ws.Range("C54").Value = 'Site locations[]'
$ws.Range("D54").Value = 'Description'
$ws.Range("E54").ClearContents()
$ws.Range("F54").ClearContents()
$ws.Range("G54").Value = 'A text description providing details about the subject. For parking changes, this describes how the proposed works affect existing car parking arrangements.'
$ws.Range("H54").Value = 'string'
$ws.Range("I54").Value = 'MAY'

# Row 55
$ws.Range("A55").ClearContents()
$ws.Range("B55").ClearContents()
$ws.Range("C55").Value = 'Site locations[]'
$ws.Range("D55").Value = 'UPRNs[]'
$ws.Range("E55").ClearContents()
$ws.Range("F55").ClearContents()
$ws.Range("G55").Value = 'Unique Property Reference Numbers (UPRNs) for properties within the site boundary'
$ws.Range("H55").Value = 'string'
$ws.Range("I55").Value = 'MAY'

# 3. Drop the now-duplicate trailing row 56 (content shifted up into row 55)
$ws.Range("A56:I56").ClearContents()
$ws.Rows.Item(56).Delete()

# 4. Re-apply merged cells for the A (top-level) and B (top-level-description) columns
$ws.Range("A2:A19").Merge()
$ws.Range("B2:B19").Merge()
$ws.Range("A20:A23").Merge()
$ws.Range("B20:B23").Merge()
$ws.Range("A24:A31").Merge()
$ws.Range("B24:B31").Merge()
$ws.Range("A32:A35").Merge()
$ws.Range("B32:B35").Merge()
$ws.Range("A36:A41").Merge()
$ws.Range("B36:B41").Merge()
$ws.Range("A42").Merge()
$ws.Range("B42").Merge()
$ws.Range("A43").Merge()
$ws.Range("B43").Merge()
$ws.Range("A44:A46").Merge()
$ws.Range("B44:B46").Merge()
$ws.Range("A47:A55").Merge()
$ws.Range("B47:B55").Merge()

# 5. Ensure the declared used range matches the new data extent
$ws.Range("A1:I55").Select()
